$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Andy (pers_personid 349)
$ws.Range("D4").Value = 194
$ws.Range("I4").Value = 804

# Row 5 - Anthony (pers_personid 350)
$ws.Range("D5").Value = 118
$ws.Range("I5").Value = 471
$ws.Range("J5").Value = 3.99
$ws.Range("K5").Value = 19.03

# Row 8 - Jon (pers_personid 357)
$ws.Range("D8").Value = 191
$ws.Range("I8").Value = 693
$ws.Range("J8").Value = 3.63

# Row 9 - Maisy (pers_personid 360)
$ws.Range("D9").Value = 118
$ws.Range("I9").Value = 434
$ws.Range("J9").Value = 3.68

# Row 10 - Mark (pers_personid 361)
$ws.Range("D10").Value = 138
$ws.Range("I10").Value = 496
$ws.Range("J10").Value = 3.59

# Row 11 - Matt (pers_personid 362)
$ws.Range("D11").Value = 186
$ws.Range("I11").Value = 697
$ws.Range("J11").Value = 3.75

# Row 12 - Pepe (pers_personid 364)
$ws.Range("D12").Value = 100
$ws.Range("I12").Value = 371
$ws.Range("J12").Value = 3.71

# Row 13 - Prashant (pers_personid 365)
$ws.Range("D13").Value = 33
$ws.Range("I13").Value = 125
$ws.Range("J13").Value = 3.79

# Row 14 - Richard (pers_personid 366)
$ws.Range("D14").Value = 137
$ws.Range("I14").Value = 585
$ws.Range("J14").Value = 4.27
